# Update the "Förändrad" (changed) date column C for all existing data rows
# (rows 2 through 338) from 45181 to 45182.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 338; $r++) {
    $ws.Cells.Item($r, 3).Value = 45182
}

# Append a new row (339) describing a new notice.
$newRow = 339
$ws.Cells.Item($newRow, 1).Value = "A 42515-2023"
$ws.Cells.Item($newRow, 2).Value = 45180
$ws.Cells.Item($newRow, 3).Value = 45182
$ws.Cells.Item($newRow, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item($newRow, 5).Value = "DOROTEA"
$ws.Cells.Item($newRow, 6).Value = "SCA"
$ws.Cells.Item($newRow, 7).Value = 27.7
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 0
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = 0
$ws.Cells.Item($newRow, 15).Value = 0
$ws.Cells.Item($newRow, 16).Value = 0
$ws.Cells.Item($newRow, 17).Value = 0
$ws.Cells.Item($newRow, 18).Value = ""

# Match the date-formatted columns (B, C) and the wrap-text "Artnamn" column
# (R) used by every other data row.
$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow, 18).WrapText = $true

# Keep row heights consistent with the rest of the sheet (15pt, explicit).
$ws.Rows.Item(338).RowHeight = 15
$ws.Rows.Item(339).RowHeight = 15
